# Updated symbol list on Fri Feb 10 10:35:02 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.84"
$ws.Range("E2").Value = "'-4.37%"
$ws.Range("D3").Value = "'40.00"
$ws.Range("E3").Value = "'-6.47%"
$ws.Range("E4").Value = "'-2.06%"
$ws.Range("D5").Value = "'0.07730"
$ws.Range("E5").Value = "'-5.64%"
$ws.Range("D6").Value = "'4.241"
$ws.Range("E6").Value = "'-1.29%"
$ws.Range("D7").Value = "'1.612"
$ws.Range("E7").Value = "'-11.36%"
$ws.Range("D8").Value = "'0.8844"
$ws.Range("E8").Value = "'-5.26%"
$ws.Range("D9").Value = "'0.09909"
$ws.Range("E9").Value = "'-10.65%"
$ws.Range("E10").Value = "'-6.60%"
$ws.Range("D11").Value = "'0.09029"
$ws.Range("E11").Value = "'-5.23%"
$ws.Range("D12").Value = "'0.04437"
$ws.Range("E12").Value = "'-5.30%"
$ws.Range("E13").Value = "'-0.39%"
$ws.Range("D14").Value = "'0.001260"
$ws.Range("E14").Value = "'-3.89%"
$ws.Range("D15").Value = "'0.005800"
$ws.Range("E15").Value = "'-1.26%"
$ws.Range("E16").Value = "'2,413.39%"
$ws.Range("D17").Value = "'3.352"
$ws.Range("E17").Value = "'-0.04%"
$ws.Range("E18").Value = "'-3.53%"
$ws.Range("E19").Value = "'-3.19%"
$ws.Range("D20").Value = "'7.097"
$ws.Range("E20").Value = "'-4.59%"
$ws.Range("E21").Value = "'-2.83%"
$ws.Range("D22").Value = "'0.2849"
$ws.Range("E22").Value = "'11.87%"
$ws.Range("D23").Value = "'0.04117"
$ws.Range("E23").Value = "'-0.87%"
$ws.Range("E24").Value = "'-4.20%"
$ws.Range("D25").Value = "'0.004079"
$ws.Range("E25").Value = "'-5.91%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'8.48%"
$ws.Range("D38").Value = "'0.02348"
$ws.Range("E38").Value = "'-14.22%"
$ws.Range("D39").Value = "'0.05210"
$ws.Range("E39").Value = "'-6.94%"
$ws.Range("D40").Value = "'0.007964"
$ws.Range("E40").Value = "'-1.02%"
$ws.Range("D41").Value = "'0.1323"
$ws.Range("E41").Value = "'-5.45%"
$ws.Range("D42").Value = "'0.006227"
$ws.Range("E42").Value = "'-4.84%"
$ws.Range("D43").Value = "'0.001953"
$ws.Range("E43").Value = "'-4.79%"
$ws.Range("D44").Value = "'0.008760"
$ws.Range("E44").Value = "'4.94%"
$ws.Range("D45").Value = "'0.3331"
$ws.Range("E45").Value = "'-5.04%"
$ws.Range("D46").Value = "'0.00006552"
$ws.Range("E46").Value = "'-5.86%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.14%"
$ws.Range("B48").Value = "'BOLO"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.003643"
$ws.Range("E48").Value = "'4.92%"
$ws.Range("B49").Value = "'CoinbaseStockToken"
$ws.Range("C49").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.007006"
$ws.Range("E49").Value = "'98.36%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.14%"
$ws.Range("E51").Value = "'0.14%"
